# The deck shipped with two independent DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Integral" design (clrScheme "Red Violet"), wired to
#                            the slide master / presentation (the active theme).
#   ppt/theme/theme2.xml -> "Office Theme" design (clrScheme "Office"), an extra
#                            theme part not referenced by any master/layout/slide.
#
# The authored edit swaps the two themes' contents (theme1 <-> theme2), so the
# deck now renders with the stock "Office Theme" palette while the former
# "Integral" colors move to the (unused) second theme part. theme2.xml is not
# reachable through the PowerPoint object model (no shape/slide/master ever
# points at it), so the only observable, COM-automatable side of this edit is
# repainting the live theme's 12-slot color scheme to the Office Theme values
# -- font scheme / format scheme are already byte-identical between the two
# themes, so the color scheme is the only thing that needs to move.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

# Target palette = the stock "Office Theme" clrScheme, in the fixed
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink slot order (RRGGBB, as in the XML).
$officeThemeColorsRGB = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $scheme.Count; $i++) {
    $rrggbb = $officeThemeColorsRGB[$i - 1]
    $r = [math]::Floor($rrggbb / 65536) % 256
    $g = [math]::Floor($rrggbb / 256) % 256
    $b = $rrggbb % 256
    # ThemeColorScheme.Item(..).RGB uses the classic Windows BGR-packed order.
    $bgr = ($b * 65536) + ($g * 256) + $r
    $scheme.Item($i).RGB = $bgr
}
